$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates per the diff (cryptos list refresh).
# Cells whose new value would otherwise be auto-coerced to a number by
# Excel (plain numeric-looking text in the "Price" column) are forced back
# to text, matching the source data which stores these as strings.

$ws.Range('D2').Value = '43.301.23'
$ws.Range('E2').Value = '  +2.70%  '
$ws.Range('D3').Value = '2.303.57'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').Value = '  -0.07%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '310.52'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +1.48%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '102.86'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +5.92%  '
$ws.Range('E7').Value = '  +1.20%  '
$ws.Range('E8').Value = '  -0.05%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.530'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +8.06%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '35.60'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +1.46%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0811'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +2.90%  '
$ws.Range('E12').Value = '  -1.11%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '6.97'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').Value = '2.660.32'
$ws.Range('E14').Value = '  +1.62%  '
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '2.289.46'
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.806'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +2.15%  '
$ws.Range('E18').Value = '  +2.72%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.25'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = '0.0₃0932'
$ws.Range('E20').Value = '  +3.14%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.17'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +2.89%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '68.01'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.37%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '240.70'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.01'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.61'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +1.16%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.03%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '3.97'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -1.83%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '24.99'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +6.38%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.30'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +8.10%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '36.58'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -2.41%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '9.62'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.39%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '169.92'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +4.69%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '5.26'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  +6.83%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '17.77'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.95%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.0740'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('E38').Value = '  -2.57%  '
$ws.Range('E39').Value = '  +3.34%  '
$ws.Range('E40').Value = '  +2.07%  '
$ws.Range('E41').Value = '  +0.98%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '4.35'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +6.84%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.30'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -1.95%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.0289'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +2.92%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '19.26'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.964.94'
$ws.Range('E46').Value = '  +0.88%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.99'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +2.60%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '9.90'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -0.30%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '55.31'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +2.82%  '
$ws.Range('E50').Value = '  +1.18%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.58'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +7.41%  '
